$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contact discovered for the remaining rows of the report (rows 30-42):
# fill column A with the new email address, matching the existing pattern
# used for the earlier contacts in the sheet (plain hyperlinked text).
$email = "eduardo.XXXXXX@hotmail.com"

for ($r = 30; $r -le 42; $r++) {
    $ws.Cells.Item($r, 1).Value = $email
}

# Row 30 gets its own single-cell hyperlink (mirrors A10 / A22 above), and
# rows 31-42 share one hyperlink range with a display string (mirrors
# A3:A9 / A11:A21 / A23:A29 above).
$ws.Hyperlinks.Add($ws.Range("A30"), "mailto:$email") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A31:A42"), "mailto:$email", $null, $null, $email) | Out-Null

# Hyperlinks.Add() re-styles the touched cells with a brand new cell
# format; reapply the workbook's existing "Hyperlink" style so these cells
# keep using the same shared style as the rest of column A.
$ws.Range("A30:A42").Style = "Hyperlink"

# Leave the selection where the user left it when they saved: the newly
# filled range, anchored at the top cell.
$ws.Range("A30:A42").Select() | Out-Null
